$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '68.507.93'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +1.55%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.919.37'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +1.40%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '484.92'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +4.18%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '147.29'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.26%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.620'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -1.94%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.999'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.03%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.722'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -3.53%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.168'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +8.10%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0000354'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +13.39%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '42.42'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -3.81%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '10.50'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +0.55%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.546.16'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +1.35%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '14.57'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -0.93%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.910.71'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +1.09%  '

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -0.47%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '19.71'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -1.70%  '

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -2.99%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '68.697.66'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +1.58%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '432.34'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.17%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '14.52'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -2.05%  '

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +1.25%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '86.85'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -2.03%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '11.42'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +11.70%  '

$ws.Range("B26").NumberFormat = "@"
$ws.Range("B26").Value = 'RenderToken'
$ws.Range("C26").NumberFormat = "@"
$ws.Range("C26").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '10.73'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +2.77%  '

$ws.Range("B27").NumberFormat = "@"
$ws.Range("B27").Value = 'PancakeSwap'
$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '3.58'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +0.48%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '38.02'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +1.03%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.89'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +6.73%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '718.04'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -3.99%  '

$ws.Range("B31").NumberFormat = "@"
$ws.Range("B31").Value = 'Hedera'
$ws.Range("C31").NumberFormat = "@"
$ws.Range("C31").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.129'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -3.70%  '

$ws.Range("B32").NumberFormat = "@"
$ws.Range("B32").Value = 'Cosmos'
$ws.Range("C32").NumberFormat = "@"
$ws.Range("C32").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '13.22'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -4.01%  '

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +2.90%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0₃0903'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +33.80%  '

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -4.63%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '58.38'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +1.75%  '

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -7.08%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.53'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -0.51%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.998'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -0.09%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.85'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +8.68%  '

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -2.22%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.05'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +11.11%  '

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +1.76%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.343'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -2.50%  '

$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = 'FirstDigitalUSD'
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.00'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -0.06%  '

$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = 'Stellar'
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.140'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -1.28%  '

$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = 'ARBITRUM'
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.16'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +1.18%  '

$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = 'LidoDAOToken'
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.40'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -1.11%  '

$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = 'Monero'
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '148.01'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +2.07%  '

$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = 'ApeXProtocol'
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '3.19'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -2.95%  '

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -2.76%  '
